$wb = $excel.ActiveWorkbook

# --- Fix OCR garbage in the "register_reason" string on the land sheet ("土地") ---
# "11rt■■貝買" -> "11rt貝買"
$landSheet = $wb.Worksheets.Item(1)
$landSheet.Cells.Item(5, 7).Value = "11rt貝買"

# --- Rename sheet 6 from "具有相當價值之財產" to "保險" (insurance) ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "保險"

# --- Rebuild sheet 6 contents: drop the unused E column and the extra 4th row,
#     and populate the remaining cells with the correct insurance data for
#     富邦人壽 (Fubon Life): two policies held by 林正二. ---
$ws.Columns.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Header-like first row (matches the pattern already used elsewhere in this sheet)
$ws.Cells.Item(1, 2).Value = "富邦人壽"
$ws.Cells.Item(1, 3).Value = "生存還本保險"
$ws.Cells.Item(1, 4).Value = "林正二"

# Row 2 (index 82)
$ws.Cells.Item(2, 1).Value = 82
$ws.Cells.Item(2, 2).Value = "富邦人壽"
$ws.Cells.Item(2, 3).Value = "生存還本保險"
$ws.Cells.Item(2, 4).Value = "林正二"

# Row 3 (index 83)
$ws.Cells.Item(3, 1).Value = 83
$ws.Cells.Item(3, 2).Value = "富邦人壽"
$ws.Cells.Item(3, 3).Value = "年金保險"
$ws.Cells.Item(3, 4).Value = "林正二"
